# Update cryptos list with fresh price/volume data and
# re-rank Hedera / EthereumClassic / NEARProtocol (rows 31-33).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.425.10"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").Value = "2.984.72"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'566.93"
$ws.Range("E5").Value = "  +2.47%  "
$ws.Range("D6").Value = "'138.55"
$ws.Range("E6").Value = "  +4.10%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("D9").Value = "2.976.86"
$ws.Range("E9").Value = "  +1.29%  "
$ws.Range("E10").Value = "  +3.85%  "
$ws.Range("D11").Value = "'5.39"
$ws.Range("E11").Value = "  +12.00%  "
$ws.Range("D12").Value = "'0.453"
$ws.Range("E12").Value = "  +1.60%  "
$ws.Range("D13").Value = "'0.0000230"
$ws.Range("E13").Value = "  +3.92%  "
$ws.Range("D14").Value = "'33.74"
$ws.Range("E14").Value = "  +2.85%  "
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "3.472.07"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "'7.05"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").Value = "2.979.88"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("D19").Value = "59.429.38"
$ws.Range("D20").Value = "'437.49"
$ws.Range("E20").Value = "  +5.00%  "
$ws.Range("D21").Value = "'13.59"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").Value = "'0.719"
$ws.Range("E22").Value = "  +3.78%  "
$ws.Range("D23").Value = "'13.38"
$ws.Range("E23").Value = "  -0.37%  "
$ws.Range("D24").Value = "'7.03"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "'79.98"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'2.21"
$ws.Range("E27").Value = "  +9.72%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").Value = "'2.55"
$ws.Range("E29").Value = "  +2.38%  "
$ws.Range("D30").Value = "'7.73"
$ws.Range("E30").Value = "  +3.85%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.106"
$ws.Range("E31").Value = "  +9.59%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'25.75"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'6.22"
$ws.Range("E33").Value = "  +4.77%  "
$ws.Range("D34").Value = "0.0₃0769"
$ws.Range("E34").Value = "  +9.44%  "
$ws.Range("D35").Value = "'5.90"
$ws.Range("E35").Value = "  +3.93%  "
$ws.Range("D36").Value = "'0.985"
$ws.Range("E36").Value = "  +4.12%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").Value = "'48.65"
$ws.Range("E38").Value = "  +0.87%  "
$ws.Range("D39").Value = "'8.65"
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("E40").Value = "  +4.40%  "
$ws.Range("D41").Value = "'401.30"
$ws.Range("E41").Value = "  +5.64%  "
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").Value = "2.746.48"
$ws.Range("E43").Value = "  +1.83%  "
$ws.Range("D44").Value = "'0.105"
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("E45").Value = "  +6.28%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'34.93"
$ws.Range("E47").Value = "  +19.27%  "
$ws.Range("D48").Value = "'122.42"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").Value = "'23.34"
$ws.Range("E51").Value = "  +1.89%  "
